$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(27, "flower/flower003.jpg", "spüren", "flower"),
    @(72, "flower/flower016.jpg", "quellen", "flower"),
    @(43, "face/face010.jpg", "leeren", "face"),
    @(109, "face/face013.jpg", "wehen", "face"),
    @(42, "flower/flower025.jpg", "holen", "flower"),
    @(70, "face/face031.jpg", "ärgern", "face"),
    @(13, "face/face003.jpg", "piepen", "face"),
    @(38, "face/face021.jpg", "hören", "face"),
    @(85, "flower/flower021.jpg", "spenden", "flower"),
    @(114, "flower/flower028.jpg", "reisen", "flower"),
    @(78, "flower/flower022.jpg", "binden", "flower"),
    @(15, "face/face011.jpg", "frischen", "face"),
    @(87, "face/face018.jpg", "grenzen", "face"),
    @(35, "face/face004.jpg", "rufen", "face"),
    @(108, "flower/flower029.jpg", "münzen", "flower"),
    @(124, "face/face008.jpg", "angeln", "face"),
    @(95, "flower/flower013.jpg", "wählen", "flower"),
    @(81, "face/face019.jpg", "kranken", "face"),
    @(63, "flower/flower026.jpg", "narren", "flower"),
    @(0, "face/face022.jpg", "prüfen", "face"),
    @(82, "flower/flower010.jpg", "streifen", "flower"),
    @(98, "face/face014.jpg", "herrschen", "face"),
    @(2, "flower/flower031.jpg", "zögern", "flower"),
    @(107, "face/face007.jpg", "proben", "face"),
    @(22, "flower/flower019.jpg", "stoppen", "flower"),
    @(67, "face/face026.jpg", "lächeln", "face"),
    @(90, "face/face005.jpg", "fließen", "face"),
    @(8, "flower/flower005.jpg", "planen", "flower"),
    @(100, "flower/flower007.jpg", "lügen", "flower"),
    @(1, "flower/flower002.jpg", "duschen", "flower"),
    @(36, "face/face009.jpg", "stoßen", "face"),
    @(34, "flower/flower001.jpg", "zielen", "flower")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $item = $data[$i]
    $ws.Cells.Item($rowNum, 2).Value = $item[0]
    $ws.Cells.Item($rowNum, 3).Value = $item[1]
    $ws.Cells.Item($rowNum, 4).Value = $item[2]
    $ws.Cells.Item($rowNum, 5).Value = $item[3]
}
